$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (exhibitions) — bump a few "want to go" counters and insert a
# new event row (南宁·AB动漫游戏嘉年华) before the existing 横州 row, shifting
# 横州 down from row 8 to row 9.
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")

$wsExpo.Range("F2").Value = 7027
$wsExpo.Range("F5").Value = 100
$wsExpo.Range("F6").Value = 1084

# Push the existing row 8 (横州·第二届海棠动漫游戏嘉年华) down to row 9.
$wsExpo.Range("A9").Value = 8
$wsExpo.Range("A8").Copy()
$wsExpo.Range("A9").PasteSpecial(-4122)

# B holds a plain "yyyy-mm-dd" looking string; force it in as text (leading
# apostrophe) so it isn't reinterpreted as a real date, then reset its
# format from the source cell so no quote-prefix/number-format lingers.
$wsExpo.Range("B9").Value = "'" + $wsExpo.Range("B8").Value2
$wsExpo.Range("B8").Copy()
$wsExpo.Range("B9").PasteSpecial(-4122)

$wsExpo.Range("C9").Value = $wsExpo.Range("C8").Value2
$wsExpo.Range("D9").Value = $wsExpo.Range("D8").Value2
$wsExpo.Range("E9").Value = $wsExpo.Range("E8").Value2
$wsExpo.Range("F9").Value = $wsExpo.Range("F8").Value2
$wsExpo.Range("G9").Value = $wsExpo.Range("G8").Value2
$wsExpo.Range("H9").Value = $wsExpo.Range("H8").Value2
$wsExpo.Range("I9").Value = $wsExpo.Range("I8").Value2
$wsExpo.Application.CutCopyMode = $false

# Write the new event into row 8 (A8/B8 keep their existing values: 7 and
# 2024-07-20, identical to what the new event needs).
$wsExpo.Range("C8").Value = "南宁·AB动漫游戏嘉年华"
$wsExpo.Range("D8").Value = "三塘南路与长虹东路交叉路口往北约50米 广西农业会展中心"
$wsExpo.Range("E8").Value = "2024.07.20 09:30-07.21 17:00"
$wsExpo.Range("F8").Value = 2
$wsExpo.Range("G8").Value = 60
$wsExpo.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=84862"
$wsExpo.Range("I8").Value = "//i1.hdslb.com/bfs/openplatform/202404/eglavDeZ1714036487217.jpeg"

# ---------------------------------------------------------------------------
# Sheet "演出" (performances) — bump the "want to go" counter for the only
# event row.
# ---------------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 9

# ---------------------------------------------------------------------------
# Sheet "全部类型" (all types) — same updates as "展览" above, plus the
# performance counter bump, since this sheet aggregates every event.
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")

$wsAll.Range("F2").Value = 7027
$wsAll.Range("F5").Value = 100
$wsAll.Range("F6").Value = 1084
$wsAll.Range("F8").Value = 9

# Push the existing row 9 (横州·第二届海棠动漫游戏嘉年华) down to row 10.
$wsAll.Range("A10").Value = 9
$wsAll.Range("A9").Copy()
$wsAll.Range("A10").PasteSpecial(-4122)

$wsAll.Range("B10").Value = "'" + $wsAll.Range("B9").Value2
$wsAll.Range("B9").Copy()
$wsAll.Range("B10").PasteSpecial(-4122)

$wsAll.Range("C10").Value = $wsAll.Range("C9").Value2
$wsAll.Range("D10").Value = $wsAll.Range("D9").Value2
$wsAll.Range("E10").Value = $wsAll.Range("E9").Value2
$wsAll.Range("F10").Value = $wsAll.Range("F9").Value2
$wsAll.Range("G10").Value = $wsAll.Range("G9").Value2
$wsAll.Range("H10").Value = $wsAll.Range("H9").Value2
$wsAll.Range("I10").Value = $wsAll.Range("I9").Value2
$wsAll.Application.CutCopyMode = $false

# Write the new event into row 9 (A9/B9 keep their existing values: 8 and
# 2024-07-20, identical to what the new event needs).
$wsAll.Range("C9").Value = "南宁·AB动漫游戏嘉年华"
$wsAll.Range("D9").Value = "三塘南路与长虹东路交叉路口往北约50米 广西农业会展中心"
$wsAll.Range("E9").Value = "2024.07.20 09:30-07.21 17:00"
$wsAll.Range("F9").Value = 2
$wsAll.Range("G9").Value = 60
$wsAll.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=84862"
$wsAll.Range("I9").Value = "//i1.hdslb.com/bfs/openplatform/202404/eglavDeZ1714036487217.jpeg"
